# Update "D score categories and intervals.xlsx":
# Add a "Notes" column (H) with a header + a wrapped notes cell, widen the
# column, increase row 2's height so the note fits, and move the selection
# to the new notes cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: width to fit the notes text.
$ws.Range("H1").ColumnWidth = 73

# Header cell H1 - plain Helvetica style (matches column default, style 1).
$ws.Range("H1").Value = "Notes"

# Notes cell H2 - wrapped text, taller row to show the whole note.
$ws.Range("H2").Value = "D score cutoffs were extracted from the IAT code employed on project implicit that gives feedback to participants. The predicted 95% CI is taken from our current results, ie the predicted 95% interval for that D score for each domain. See analysis_iat_D.html and analysis_iat_D.Rmd"
$ws.Range("H2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 68

# Move the active selection to the new notes cell, as in the saved file.
[void]$ws.Range("H3").Select()
